# Add new test results: overwrite the first measurement block (A1:J5)
# with the latest run's readings. The difference rows (18:22) recompute
# automatically since they hold formulas referencing A1:J5 and A9:J13.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 55.6
$ws.Range("B1").Value = 66.6
$ws.Range("C1").Value = 71.9
$ws.Range("D1").Value = 82.2
$ws.Range("E1").Value = 72
$ws.Range("F1").Value = 70.7
$ws.Range("G1").Value = 84.6
$ws.Range("H1").Value = 71.3
$ws.Range("I1").Value = 123.6
$ws.Range("J1").Value = 148.2

$ws.Range("A2").Value = 70.6
$ws.Range("B2").Value = 69.7
$ws.Range("C2").Value = 147.4
$ws.Range("D2").Value = 101.3
$ws.Range("E2").Value = 114.9
$ws.Range("F2").Value = 85.6
$ws.Range("G2").Value = 66.6
$ws.Range("H2").Value = 181.4
$ws.Range("I2").Value = 137
$ws.Range("J2").Value = 128.1

$ws.Range("A3").Value = 85.5
$ws.Range("B3").Value = 54.3
$ws.Range("C3").Value = 145.5
$ws.Range("D3").Value = 148.5
$ws.Range("E3").Value = 111.1
$ws.Range("F3").Value = 77.2
$ws.Range("G3").Value = 60.5
$ws.Range("H3").Value = 159.4
$ws.Range("I3").Value = 136.1
$ws.Range("J3").Value = 148.1

$ws.Range("A4").Value = 72.1
$ws.Range("B4").Value = 53.3
$ws.Range("C4").Value = 121.7
$ws.Range("D4").Value = 116.2
$ws.Range("E4").Value = 149.3
$ws.Range("F4").Value = 63.5
$ws.Range("G4").Value = 47.3
$ws.Range("H4").Value = 151.9
$ws.Range("I4").Value = 161.8
$ws.Range("J4").Value = 171.2

$ws.Range("A5").Value = 73.1
$ws.Range("B5").Value = 57.9
$ws.Range("C5").Value = 109.5
$ws.Range("D5").Value = 141.4
$ws.Range("E5").Value = 129
$ws.Range("F5").Value = 78.5
$ws.Range("G5").Value = 78.4
$ws.Range("H5").Value = 187
$ws.Range("I5").Value = 140.1
$ws.Range("J5").Value = 131.1
